# Daily attendance processing - reorder "Recorded By" email lists in column G
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$values = @{
    'G2' = 'rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
    'G4' = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg'
    'G5' = 'Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
    'G6' = 'maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
    'G7' = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
    'G8' = 'Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg'
    'G9' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G10' = 'Alshimaa_khaled@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg'
    'G12' = 'mohamed.saleem@med.asu.edu.eg, System'
    'G14' = 'marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg'
    'G15' = 'nourhanmohamed@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
    'G16' = 'marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, randa.rabea@med.asu.edu.eg'
    'G17' = 'nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
    'G18' = 'Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
    'G19' = 'rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
    'G21' = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg'
    'G22' = 'Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
    'G23' = 'maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
    'G24' = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
    'G25' = 'Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg'
    'G26' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G27' = 'Alshimaa_khaled@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg'
    'G30' = 'Sarah.Mahdy@med.asu.edu.eg, mariam.g.sanad@med.asu.edu.eg'
    'G31' = 'marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg'
    'G32' = 'nourhanmohamed@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
    'G33' = 'marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, randa.rabea@med.asu.edu.eg'
    'G34' = 'nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg'
    'G35' = 'abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
    'G36' = 'rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
    'G37' = 'NadaMohamed@med.asu.edu.eg, System, Kerelos.zareef@med.asu.edu.eg'
    'G39' = 'Omnia.Mohammed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg'
    'G40' = 'maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
    'G43' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G44' = 'Alshimaa_khaled@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg'
    'G45' = 'mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com'
    'G46' = 'mohamed.saleem@med.asu.edu.eg, System'
    'G48' = 'marian.samir@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
    'G50' = 'nancy.abdelshafy@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg'
    'G51' = 'abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
    'G52' = 'abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
    'G53' = 'rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
    'G54' = 'NadaMohamed@med.asu.edu.eg, System, Kerelos.zareef@med.asu.edu.eg'
    'G56' = 'Omnia.Mohammed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg'
    'G57' = 'maryam.ahmed@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, maimustafa@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
    'G60' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G61' = 'Alshimaa_khaled@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg'
    'G62' = 'mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com'
    'G65' = 'marian.samir@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg'
    'G67' = 'nancy.abdelshafy@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg'
    'G68' = 'abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
    'G69' = 'abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg'
    'G70' = 'rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
    'G72' = 'Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
    'G73' = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg'
    'G74' = 'Sara_nabil@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
    'G75' = 'maryam.ahmed@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
    'G76' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
    'G77' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G78' = 'Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G79' = 'Alshimaa_khaled@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg'
    'G84' = 'marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg'
    'G85' = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nahla91@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
    'G86' = 'Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
    'G87' = 'rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
    'G88' = 'NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
    'G89' = 'aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg'
    'G90' = 'Sara_nabil@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
    'G91' = 'maryam.ahmed@med.asu.edu.eg, Eman.m.elsaid@med.asu.edu.eg, merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg'
    'G92' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G93' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G94' = 'Alshimaa_khaled@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg'
    'G98' = 'nourhanmohamed@med.asu.edu.eg, user@user.com, afaf.abdallah@med.asu.edu.eg, nourhanhosni@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg'
    'G99' = 'user@user.com, Walaa.h.ghanima@med.asu.edu.eg'
    'G101' = 'marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg'
    'G102' = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, nahla91@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
    'G103' = 'abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, gigi@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
    'G104' = 'rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
    'G106' = 'Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
    'G107' = 'Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
    'G108' = 'merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
    'G109' = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
    'G110' = 'Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
    'G111' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G112' = 'Alshimaa_khaled@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg'
    'G114' = 'mohamed.saleem@med.asu.edu.eg, System'
    'G115' = 'Sarah.Mahdy@med.asu.edu.eg, mariam.g.sanad@med.asu.edu.eg'
    'G116' = 'nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg'
    'G117' = 'Amr-Saeed@med.asu.edu.eg, user@user.com, takyas@med.asu.edu.eg, enas.omran@med.asu.edu.eg'
    'G118' = 'marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg'
    'G119' = 'nardine.alfonse@med.asu.edu.eg, marina_atef@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nahla91@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
    'G120' = 'abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, gigi@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
    'G121' = 'rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
    'G123' = 'Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg'
    'G124' = 'Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
    'G126' = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
    'G127' = 'Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
    'G128' = 'Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G129' = 'Alshimaa_khaled@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg'
    'G131' = 'mohamed.saleem@med.asu.edu.eg, System'
    'G132' = 'Sarah.Mahdy@med.asu.edu.eg, mariam.g.sanad@med.asu.edu.eg'
    'G133' = 'nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg'
    'G134' = 'Amr-Saeed@med.asu.edu.eg, user@user.com, takyas@med.asu.edu.eg, enas.omran@med.asu.edu.eg'
    'G135' = 'marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg'
    'G136' = 'nardine.alfonse@med.asu.edu.eg, marina_atef@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nahla91@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
    'G137' = 'abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, gigi@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
    'G138' = 'rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg'
    'G142' = 'merna.said@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, basma.hamed@med.asu.edu.eg'
    'G143' = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
    'G144' = 'mohamed.saleem@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G145' = 'Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G146' = 'mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg'
    'G148' = 'Sarah.Mahdy@med.asu.edu.eg, mariam.g.sanad@med.asu.edu.eg'
    'G149' = 'user@user.com, Walaa.h.ghanima@med.asu.edu.eg'
    'G151' = 'nourhanmohamed@med.asu.edu.eg, hana.amr@med.asu.edu.eg, System, randa.rabea@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
    'G152' = 'marian.samir@med.asu.edu.eg, ManaratAleslam-Zefan@med.asu.edu.eg'
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}

